$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.925.03"
$ws.Range("E2").Value = "  -1.07%  "
$ws.Range("D3").Value = "2.193.90"
$ws.Range("E3").Value = "  -2.27%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.68"
$ws.Range("E5").Value = "  -4.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "89.25"
$ws.Range("E6").Value = "  -5.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.563"
$ws.Range("E7").Value = "  -1.49%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.484"
$ws.Range("E9").Value = "  -7.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.07"
$ws.Range("E10").Value = "  -8.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0773"
$ws.Range("E11").Value = "  -4.63%  "
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.82"
$ws.Range("E13").Value = "  -5.37%  "
$ws.Range("D14").Value = "2.528.24"
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("D15").Value = "2.262.15"
$ws.Range("E15").Value = "  -4.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.16"
$ws.Range("E16").Value = "  -3.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.774"
$ws.Range("E17").Value = "  -7.79%  "
$ws.Range("D18").Value = "43.626.57"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("D19").Value = "0.0₃0891"
$ws.Range("E19").Value = "  -7.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.83"
$ws.Range("E20").Value = "  -9.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.83"
$ws.Range("E21").Value = "  -12.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "63.20"
$ws.Range("E22").Value = "  -3.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.80"
$ws.Range("E23").Value = "  -2.41%  "
$ws.Range("E24").Value = "  -9.10%  "
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  -8.77%  "
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.42"
$ws.Range("E28").Value = "  -5.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.26"
$ws.Range("E29").Value = "  -6.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.23"
$ws.Range("E30").Value = "  -4.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "148.30"
$ws.Range("E31").Value = "  -3.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.28"
$ws.Range("E32").Value = "  -11.27%  "
$ws.Range("E33").Value = "  -5.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0736"
$ws.Range("E34").Value = "  -7.94%  "
$ws.Range("E35").Value = "  -4.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.83"
$ws.Range("E36").Value = "  -8.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.103"
$ws.Range("E37").Value = "  -6.03%  "
$ws.Range("E38").Value = "  -8.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0285"
$ws.Range("E39").Value = "  -6.10%  "
$ws.Range("E40").Value = "  -8.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.07"
$ws.Range("E41").Value = "  -11.98%  "
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.02"
$ws.Range("E43").Value = "  -11.19%  "
$ws.Range("D44").Value = "1.792.75"
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("E45").Value = "  +3.59%  "
$ws.Range("E46").Value = "  +11.43%  "
$ws.Range("E47").Value = "  -10.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "72.30"
$ws.Range("E48").Value = "  -10.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "91.80"
$ws.Range("E49").Value = "  -8.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.98"
$ws.Range("E50").Value = "  -8.44%  "
$ws.Range("D51").Value = "2.411.01"
$ws.Range("E51").Value = "  -2.25%  "
